$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39, pushing the existing row 39..78 data down to 40..79.
$ws.Rows("39:39").Insert()

# Fill in the new row 39 with fresh data.
$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 44705
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = 100112031
$ws.Range("G39").Value = "Poroto verde"
$ws.Range("H39").Value = "Magnum"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 35
$ws.Range("K39").Value = 30000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = 30000
$ws.Range("N39").Value = "`$/malla 25 kilos"
$ws.Range("O39").Value = "Perú"
$ws.Range("P39").Value = 1200
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
